$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H95").Value = 359967.75
$ws.Range("J95").Value = 359967.75
$ws.Range("L95").Value = 359967.75
$ws.Range("N95").Value = -365459.75
$ws.Range("H98").Value = 534655.9
$ws.Range("I98").Value = 933437.0600000001
$ws.Range("J98").Value = 2947.5557
$ws.Range("K98").Value = 933437.0600000001
$ws.Range("L98").Value = 2947.5557
$ws.Range("M98").Value = -931939.0600000001
$ws.Range("N98").Value = -5943.5557
$ws.Range("H122").Value = 534655.9
$ws.Range("I122").Value = 933437.0600000001
$ws.Range("J122").Value = 2947.5557
$ws.Range("K122").Value = 2800311.18
$ws.Range("L122").Value = 8842.667099999999
$ws.Range("M122").Value = -2797861.18
$ws.Range("N122").Value = -13742.6671
$ws.Range("H126").Value = 36000
$ws.Range("J126").Value = 36000
$ws.Range("L126").Value = 36000
$ws.Range("N126").Value = -45880
$ws.Range("H128").Value = 80780
$ws.Range("J128").Value = 80780
$ws.Range("L128").Value = 80780
$ws.Range("N128").Value = -90740
$ws.Range("H133").Value = 24028.334
$ws.Range("J133").Value = 24028.334
$ws.Range("L133").Value = 24028.334
$ws.Range("N133").Value = -34148.334
$ws.Range("H136").Value = 48333.332
$ws.Range("J136").Value = 48333.332
$ws.Range("L136").Value = 48333.332
$ws.Range("N136").Value = -58533.332
$ws.Range("H137").Value = 29413090
$ws.Range("I137").Value = 43479224
$ws.Range("J137").Value = 2088.9092
$ws.Range("K137").Value = 130437672
$ws.Range("L137").Value = 6266.7276
$ws.Range("M137").Value = -130435122
$ws.Range("N137").Value = -11366.7276
$ws.Range("H139").Value = 50000
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 50000
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 50000
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -60280
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17752.703
$ws.Range("I32").Value = 2065.4575
$ws.Range("J32").Value = 202862.2
$ws.Range("K32").Value = 2065.4575
$ws.Range("L32").Value = 202862.2
$ws.Range("M32").Value = -1778.4575
$ws.Range("N32").Value = -203436.2
$ws.Range("H110").Value = 744.4
$ws.Range("I110").Value = 630
$ws.Range("J110").Value = 973.2
$ws.Range("K110").Value = 630
$ws.Range("L110").Value = 973.2
$ws.Range("M110").Value = 1415
$ws.Range("N110").Value = -5063.2
$ws.Range("H133").Value = 43890.25
$ws.Range("J133").Value = 43890.25
$ws.Range("L133").Value = 43890.25
$ws.Range("N133").Value = -48950.25
$ws.Range("H138").Value = 55250
$ws.Range("J138").Value = 55250
$ws.Range("L138").Value = 55250
$ws.Range("N138").Value = -65530
$ws.Range("H139").Value = 41571.668
$ws.Range("J139").Value = 41571.668
$ws.Range("L139").Value = 41571.668
$ws.Range("N139").Value = -51851.668
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1046.9546
$ws.Range("I94").Value = 897.8823
$ws.Range("J94").Value = 1553.8
$ws.Range("K94").Value = 897.8823
$ws.Range("L94").Value = 1553.8
$ws.Range("M94").Value = -446.8823
$ws.Range("N94").Value = -2455.8
$ws.Range("H134").Value = 1876.7123
$ws.Range("I134").Value = 1101.3148
$ws.Range("J134").Value = 4080.4736
$ws.Range("K134").Value = 3303.9444
$ws.Range("L134").Value = 12241.4208
$ws.Range("M134").Value = -768.9444000000003
$ws.Range("N134").Value = -17311.4208
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2568.2778
$ws.Range("I31").Value = 1358.5714
$ws.Range("J31").Value = 3338.0908
$ws.Range("K31").Value = 1358.5714
$ws.Range("L31").Value = 3338.0908
$ws.Range("M31").Value = -1063.5714
$ws.Range("N31").Value = -3928.0908
$ws.Range("H34").Value = 2568.2778
$ws.Range("I34").Value = 1358.5714
$ws.Range("J34").Value = 3338.0908
$ws.Range("K34").Value = 1358.5714
$ws.Range("L34").Value = 3338.0908
$ws.Range("M34").Value = -1156.5714
$ws.Range("N34").Value = -3742.0908
$ws.Range("H58").Value = 2166.634
$ws.Range("I58").Value = 864
$ws.Range("K58").Value = 864
$ws.Range("M58").Value = -661
$ws.Range("H132").Value = 1969.9678
$ws.Range("I132").Value = 1422.76
$ws.Range("J132").Value = 4250
$ws.Range("K132").Value = 4268.28
$ws.Range("L132").Value = 12750
$ws.Range("M132").Value = -1738.28
$ws.Range("N132").Value = -17810
$ws.Range("H136").Value = 2166.634
$ws.Range("I136").Value = 864
$ws.Range("K136").Value = 2592
$ws.Range("M136").Value = -42
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 10001
$ws.Range("J95").Value = 10001
$ws.Range("L95").Value = 10001
$ws.Range("N95").Value = -15493
$ws.Range("H102").Value = 3938.1052
$ws.Range("I102").Value = 3556.7273
$ws.Range("J102").Value = 4462.5
$ws.Range("K102").Value = 3556.7273
$ws.Range("L102").Value = 4462.5
$ws.Range("M102").Value = -1934.7273
$ws.Range("N102").Value = -7706.5
$ws.Range("H122").Value = 586070.9
$ws.Range("I122").Value = 855665.1
$ws.Range("J122").Value = 1950
$ws.Range("K122").Value = 2566995.3
$ws.Range("L122").Value = 5850
$ws.Range("M122").Value = -2564545.3
$ws.Range("N122").Value = -10750
$ws.Range("H126").Value = 2548.4707
$ws.Range("I126").Value = 2700
$ws.Range("J126").Value = 2528.2666
$ws.Range("K126").Value = 8100
$ws.Range("L126").Value = 7584.7998
$ws.Range("M126").Value = -5630
$ws.Range("N126").Value = -12524.7998
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3240.2
$ws.Range("I7").Value = 2660.8
$ws.Range("K7").Value = 2660.8
$ws.Range("M7").Value = -2548.8
$ws.Range("H40").Value = 3396.9119
$ws.Range("I40").Value = 2417.647
$ws.Range("J40").Value = 4376.1763
$ws.Range("K40").Value = 2417.647
$ws.Range("L40").Value = 4376.1763
$ws.Range("M40").Value = -2281.647
$ws.Range("N40").Value = -4648.1763
$ws.Range("H61").Value = 6647.9653
$ws.Range("I61").Value = 7061.7144
$ws.Range("J61").Value = 5561.875
$ws.Range("K61").Value = 7061.7144
$ws.Range("L61").Value = 5561.875
$ws.Range("M61").Value = -6859.7144
$ws.Range("N61").Value = -5965.875
$ws.Range("H113").Value = 6647.9653
$ws.Range("I113").Value = 7061.7144
$ws.Range("J113").Value = 5561.875
$ws.Range("K113").Value = 7061.7144
$ws.Range("L113").Value = 5561.875
$ws.Range("M113").Value = -4891.7144
$ws.Range("N113").Value = -9901.875
$ws.Range("H122").Value = 3171.7778
$ws.Range("I122").Value = 1980.375
$ws.Range("J122").Value = 3673.4211
$ws.Range("K122").Value = 5941.125
$ws.Range("L122").Value = 11020.2633
$ws.Range("M122").Value = -3491.125
$ws.Range("N122").Value = -15920.2633
$ws.Range("H126").Value = 3240.2
$ws.Range("I126").Value = 2660.8
$ws.Range("K126").Value = 7982.400000000001
$ws.Range("M126").Value = -5512.400000000001
$ws.Range("H132").Value = 3407.4878
$ws.Range("I132").Value = 2364.6667
$ws.Range("J132").Value = 5418.643
$ws.Range("K132").Value = 7094.000100000001
$ws.Range("L132").Value = 16255.929
$ws.Range("M132").Value = -4564.000100000001
$ws.Range("N132").Value = -21315.929
$ws.Range("H136").Value = 4687.2583
$ws.Range("I136").Value = 2900.7083
$ws.Range("J136").Value = 10812.571
$ws.Range("K136").Value = 8702.124899999999
$ws.Range("L136").Value = 32437.713
$ws.Range("M136").Value = -6152.124899999999
$ws.Range("N136").Value = -37537.713
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 78475.16
$ws.Range("I122").Value = 126059.625
$ws.Range("J122").Value = 2340
$ws.Range("K122").Value = 378178.875
$ws.Range("L122").Value = 7020
$ws.Range("M122").Value = -375728.875
$ws.Range("N122").Value = -11920
$ws.Range("H126").Value = 126363.25
$ws.Range("I126").Value = 333767
$ws.Range("K126").Value = 1001301
$ws.Range("M126").Value = -998831
$ws.Range("H132").Value = 10206924
$ws.Range("I132").Value = 12502783
$ws.Range("J132").Value = 3107.5557
$ws.Range("K132").Value = 37508349
$ws.Range("L132").Value = 9322.667099999999
$ws.Range("M132").Value = -37505819
$ws.Range("N132").Value = -14382.6671
$ws.Range("H136").Value = 9553264
$ws.Range("I136").Value = 9553264
$ws.Range("K136").Value = 28659792
$ws.Range("M136").Value = -28657242
$ws.Range("H140").Value = 75071.5
$ws.Range("J140").Value = 75071.5
$ws.Range("L140").Value = 75071.5
$ws.Range("N140").Value = -85431.5
